# Updates crypto price/volume data to the latest values scraped on Mon Jan  2 09:29:05 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'247.37"
$ws.Range("E2").Value = "'1.48%"
$ws.Range("D3").Value = "'30.22"
$ws.Range("E3").Value = "'12.16%"
$ws.Range("D4").Value = "'5.175"
$ws.Range("E4").Value = "'0.31%"
$ws.Range("D5").Value = "'0.05747"
$ws.Range("E5").Value = "'2.23%"
$ws.Range("D6").Value = "'6.587"
$ws.Range("E6").Value = "'1.59%"
$ws.Range("D7").Value = "'0.8594"
$ws.Range("E7").Value = "'5.42%"
$ws.Range("E8").Value = "'5.21%"
$ws.Range("D9").Value = "'0.1360"
$ws.Range("E9").Value = "'2.53%"
$ws.Range("D11").Value = "'0.02924"
$ws.Range("E11").Value = "'1.26%"
$ws.Range("D12").Value = "'0.09394"
$ws.Range("E12").Value = "'0.14%"
$ws.Range("D13").Value = "'0.001510"
$ws.Range("E13").Value = "'-0.30%"
$ws.Range("D14").Value = "'0.04127"
$ws.Range("E14").Value = "'-8.70%"
$ws.Range("D15").Value = "'0.0005991"
$ws.Range("E15").Value = "'0.07%"
$ws.Range("D16").Value = "'0.006130"
$ws.Range("E16").Value = "'-0.60%"
$ws.Range("D17").Value = "'3.504"
$ws.Range("E17").Value = "'-2.93%"
$ws.Range("D18").Value = "'3.040"
$ws.Range("E18").Value = "'0.58%"
$ws.Range("D19").Value = "'2.269"
$ws.Range("E19").Value = "'-1.65%"
$ws.Range("E20").Value = "'2.34%"
$ws.Range("D21").Value = "'0.03300"
$ws.Range("E21").Value = "'6.37%"
$ws.Range("D22").Value = "'0.1306"
$ws.Range("E22").Value = "'1.17%"
$ws.Range("D23").Value = "'3.612"
$ws.Range("E23").Value = "'-3.35%"
$ws.Range("D24").Value = "'0.1378"
$ws.Range("E24").Value = "'0.32%"
$ws.Range("E25").Value = "'-0.96%"
$ws.Range("D26").Value = "'0.004508"
$ws.Range("E26").Value = "'0.40%"
$ws.Range("D27").Value = "'0.0001179"
$ws.Range("E27").Value = "'20.37%"
$ws.Range("D28").Value = "'0.0001389"
$ws.Range("E28").Value = "'-0.67%"
$ws.Range("D40").Value = "'0.03787"
$ws.Range("E40").Value = "'4.16%"
$ws.Range("D41").Value = "'0.005776"
$ws.Range("E41").Value = "'-4.49%"
$ws.Range("E42").Value = "'1.87%"
$ws.Range("D43").Value = "'0.002298"
$ws.Range("E43").Value = "'-11.56%"
$ws.Range("D44").Value = "'0.01017"
$ws.Range("E44").Value = "'24.37%"
$ws.Range("D45").Value = "'0.00005082"
$ws.Range("E45").Value = "'-4.46%"
$ws.Range("E46").Value = "'-0.03%"
$ws.Range("D47").Value = "'0.08887"
$ws.Range("E47").Value = "'-18.43%"
$ws.Range("D48").Value = "'0.002745"
$ws.Range("E48").Value = "'10.57%"
$ws.Range("E49").Value = "'-0.03%"
$ws.Range("E50").Value = "'-0.03%"